$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.101.86"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.824.11"
$ws.Range("E3").Value = "  -0.65%  "

# Row 4
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").Value = "'241.34"
$ws.Range("E5").Value = "  -1.22%  "

# Row 6
$ws.Range("D6").Value = "'0.6143"
$ws.Range("E6").Value = "  -2.19%  "

# Row 7
$ws.Range("D7").Value = "'0.9986"
$ws.Range("E7").Value = "  -0.40%  "

# Row 8
$ws.Range("D8").Value = "'0.07326"
$ws.Range("E8").Value = "  -2.62%  "

# Row 9
$ws.Range("D9").Value = "'0.2886"
$ws.Range("E9").Value = "  -1.38%  "

# Row 10
$ws.Range("D10").Value = "'22.89"
$ws.Range("E10").Value = "  -1.38%  "

# Row 11
$ws.Range("E11").Value = "  -0.72%  "

# Row 12
$ws.Range("D12").Value = "1.808.18"
$ws.Range("E12").Value = "  -1.38%  "

# Row 13
$ws.Range("D13").Value = "'4.950"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").Value = "'0.6610"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15
$ws.Range("D15").Value = "'81.61"
$ws.Range("E15").Value = "  -1.30%  "

# Row 16
$ws.Range("D16").Value = "'0.000008946"
$ws.Range("E16").Value = "  -4.63%  "

# Row 17
$ws.Range("D17").Value = "'5.860"
$ws.Range("E17").Value = "  -2.27%  "

# Row 18
$ws.Range("D18").Value = "29.047.79"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$ws.Range("D19").Value = "2.042.31"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("D20").Value = "'236.02"
$ws.Range("E20").Value = "  +5.47%  "

# Row 21
$ws.Range("E21").Value = "  -1.17%  "

# Row 22
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.64%  "

# Row 23
$ws.Range("D23").Value = "'7.144"
$ws.Range("E23").Value = "  +0.45%  "

# Row 24
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.26%  "

# Row 25
$ws.Range("D25").Value = "'158.49"
$ws.Range("E25").Value = "  -1.05%  "

# Row 26
$ws.Range("D26").Value = "'0.1410"
$ws.Range("E26").Value = "  +0.96%  "

# Row 27
$ws.Range("D27").Value = "'8.428"

# Row 28
$ws.Range("D28").Value = "'17.62"
$ws.Range("E28").Value = "  -1.76%  "

# Row 29
$ws.Range("D29").Value = "'1.481"
$ws.Range("E29").Value = "  -1.27%  "

# Row 30
$ws.Range("D30").Value = "'0.05582"
$ws.Range("E30").Value = "  -1.57%  "

# Row 31
$ws.Range("D31").Value = "'4.087"
$ws.Range("E31").Value = "  +0.47%  "

# Row 32
$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  -1.40%  "

# Row 33
$ws.Range("D33").Value = "'1.204"
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("D34").Value = "'1.824"

# Row 35
$ws.Range("D35").Value = "'0.7336"
$ws.Range("E35").Value = "  -1.83%  "

# Row 36
$ws.Range("E36").Value = "  -0.90%  "

# Row 37
$ws.Range("D37").Value = "'2.621"
$ws.Range("E37").Value = "  -1.97%  "

# Row 38
$ws.Range("D38").Value = "'2.822"
$ws.Range("E38").Value = "  +2.06%  "

# Row 39
$ws.Range("D39").Value = "1.201.87"
$ws.Range("E39").Value = "  -1.64%  "

# Row 40
$ws.Range("D40").Value = "'0.01755"
$ws.Range("E40").Value = "  -1.67%  "

# Row 41
$ws.Range("D41").Value = "'6.393"
$ws.Range("E41").Value = "  -1.90%  "

# Row 42
$ws.Range("D42").Value = "'0.8922"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43
$ws.Range("D43").Value = "'0.9983"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").Value = "'100.60"
$ws.Range("E44").Value = "  -1.51%  "

# Row 45
$ws.Range("D45").Value = "1.958.00"
$ws.Range("E45").Value = "  -0.94%  "

# Row 46
$ws.Range("D46").Value = "'64.61"

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000121"
$ws.Range("E47").Value = "  -4.74%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5078"
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("D49").Value = "'9.080"
$ws.Range("E49").Value = "  +0.85%  "

# Row 50
$ws.Range("D50").Value = "'0.3988"
$ws.Range("E50").Value = "  -2.31%  "

# Row 51
$ws.Range("D51").Value = "'0.05793"
